# Auto-generated edit script: updates cryptos list values per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "51.821.18"
$ws.Cells.Item(2, 5).Value = "  +0.19%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.986.56"
$ws.Cells.Item(3, 5).Value = "  +1.65%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.18%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "378.51"
$ws.Cells.Item(5, 5).Value = "  +7.20%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "105.18"
$ws.Cells.Item(6, 5).Value = "  -0.17%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.548"
$ws.Cells.Item(7, 5).Value = "  -0.20%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.16%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.602"
$ws.Cells.Item(9, 5).Value = "  +0.81%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "37.70"
$ws.Cells.Item(10, 5).Value = "  +0.74%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -0.33%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0848"
$ws.Cells.Item(12, 5).Value = "  +0.23%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "18.75"
$ws.Cells.Item(13, 5).Value = "  -0.10%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "3.445.30"
$ws.Cells.Item(14, 5).Value = "  +1.18%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "7.50"
$ws.Cells.Item(15, 5).Value = "  +0.82%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "2.968.46"
$ws.Cells.Item(16, 5).Value = "  +1.12%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.964"
$ws.Cells.Item(17, 5).Value = "  -1.32%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "51.863.66"
$ws.Cells.Item(18, 5).Value = "  +0.41%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "3.49"
$ws.Cells.Item(19, 5).Value = "  +5.03%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "7.45"
$ws.Cells.Item(20, 5).Value = "  +2.04%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "13.24"
$ws.Cells.Item(21, 5).Value = "  +0.14%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "0.0₃0963"
$ws.Cells.Item(22, 5).Value = "  +0.74%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "68.90"
$ws.Cells.Item(23, 5).Value = "  -0.03%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "264.09"
$ws.Cells.Item(24, 5).Value = "  -0.64%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +4.56%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "7.48"
$ws.Cells.Item(26, 5).Value = "  +19.88%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.171"
$ws.Cells.Item(27, 5).Value = "  -2.55%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "4.17"
$ws.Cells.Item(28, 5).Value = "  -3.80%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "7.49"
$ws.Cells.Item(29, 5).Value = "  +4.17%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  -0.02%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "26.15"
$ws.Cells.Item(31, 5).Value = "  -1.24%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.105"
$ws.Cells.Item(32, 5).Value = "  -2.86%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "9.97"
$ws.Cells.Item(33, 5).Value = "  -0.58%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "34.76"
$ws.Cells.Item(34, 5).Value = "  -2.22%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "OKB"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "51.74"
$ws.Cells.Item(35, 5).Value = "  +1.69%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.08"
$ws.Cells.Item(36, 5).Value = "  -4.03%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.0439"
$ws.Cells.Item(37, 5).Value = "  +2.93%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +0.33%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "3.09"
$ws.Cells.Item(39, 5).Value = "  -3.97%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "17.49"
$ws.Cells.Item(40, 5).Value = "  +1.30%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.68"
$ws.Cells.Item(41, 5).Value = "  -6.07%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.87"
$ws.Cells.Item(42, 5).Value = "  -1.95%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  +0.96%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "124.10"
$ws.Cells.Item(44, 5).Value = "  +2.82%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "22.12"
$ws.Cells.Item(45, 5).Value = "  -3.87%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.283"
$ws.Cells.Item(46, 5).Value = "  +19.13%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -3.10%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "2.041.79"
$ws.Cells.Item(48, 5).Value = "  -2.60%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "NEARProtocol"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "3.25"
$ws.Cells.Item(49, 5).Value = "  +0.32%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "ApeXProtocol"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(50, 5).Value = "  +0.94%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0336"
$ws.Cells.Item(51, 5).Value = "  +5.50%  "
